# Apply activity-log updates to the "Activity Log" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70: team review entry with its start/end time (set before row 69's text
# so new shared strings land in the same order as the authoritative edit).
$ws.Range("D70").Value = 0.98402777777777783
$ws.Range("E70").Value = 0.99097222222222225
$ws.Range("G70").Value = "Reviewed report together with team member for possible issues."

# Row 69: clarify that the second Canvas submission was completed.
$ws.Range("G69").Value = "Submitted second attempt to Canvas. DONE"

# Row 71: revision entry with its start/end time.
$ws.Range("D71").Value = 0.99097222222222225
$ws.Range("E71").Value = 0.99652777777777779
$ws.Range("G71").Value = "Revised report together with team member for clarity"

# Row 72: final submission entry with its start/end time.
$ws.Range("D72").Value = 0.99652777777777779
$ws.Range("E72").Value = 0.99930555555555556
$ws.Range("G72").Value = "Last attempt of submission to Canvas. DONE"

# Update the sheet's active selection to match the saved view state.
$ws.Range("F69").Select() | Out-Null
